# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.215.95"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.644.11"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'596.49"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "'156.33"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +5.25%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").Value = "'5.25"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "'28.09"
$ws.Range("D14").Value = "'0.0000191"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "3.127.55"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "68.300.84"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "2.654.67"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'11.37"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "'362.52"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "'7.44"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  +3.38%  "
$ws.Range("D22").Value = "'4.82"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'2.07"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'75.09"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'9.74"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").Value = "2.790.30"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'557.89"
$ws.Range("E30").Value = "  -3.60%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").Value = "'161.12"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").Value = "'19.44"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").Value = "0.0₆0340"
$ws.Range("E42").Value = "  +5.55%  "
$ws.Range("D43").Value = "'17.79"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "'2.62"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'40.35"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "'158.55"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").Value = "'3.74"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'21.97"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "'0.0785"
$ws.Range("E51").Value = "  +0.62%  "
